$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "Segment - Red"
$ws.Range("E6").Value = "Segment - Red"
$ws.Range("E7").Value = "Segment - Red"
$ws.Range("E10").Value = "Segment - Green"
$ws.Range("E23").Value = "Segment - Red"
$ws.Range("E24").Value = "Segment - Red"
